$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Add new "Must Read" column (C) with header + Yes/No values for the two data rows.
$ws.Range("C1").Value() = "Must Read"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C2").Value() = "Yes"
$ws.Range("C3").Value() = "No"

# Move the active selection to C4, matching the post-edit workbook state.
$ws.Range("C4").Select()
